# counts_2or_fold5.xlsx update:
# "ScreenRecStarted" category is renamed to "0_unstated" throughout the
# transition-matrix sheet (header cell + the four "From ScreenRecStarted"
# row labels), and the active selection moves to G12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header (row 1) for the "ScreenRecStarted" destination bucket.
$ws.Range("G1").Value = "0_unstated"

# Row labels (column A) for the "ScreenRecStarted" source bucket.
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Move/restore the active cell selection to G12.
$ws.Range("G12").Select() | Out-Null
